$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.117.41'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.349.14'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.31%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '545.47'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.15'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.22%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.346.56'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.23%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.51'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.36%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.92%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.334'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.86'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.766.60'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '60.087.21'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.365.14'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.69'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.15'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.80'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +6.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '313.82'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.40'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.59%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.67%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +6.81%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '171.59'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.64%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +11.60%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.93'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.78%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +15.16%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.04'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.03%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.15'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +7.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '318.97'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +11.87%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '141.64'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.65%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.74%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.39'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +7.90%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.05%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.61%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0213'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +19.23%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.89%  '
